$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Go back and ensure the class diagram is up to date." paragraph:
#    it was highlighted green + struck-through; the highlight is being
#    removed (the green highlight moves to the "huge comment/class diagram
#    sweep" item below, in step 7) while the strike-through stays / becomes
#    paragraph-level.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Go back and ensure the class diagram*") {
        $p.Range.Font.HighlightColorIndex = 0   # wdNoHighlight
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 2) Relocate the hidden "_GoBack" bookmark from the "Refactor the existing
#    Unit code..." paragraph to the end of the "Create the Floor2.cs" one
#    (right after the run's text, collapsed/zero-length).
#    A plain collapsed Range exactly at end-of-paragraph-text trips an
#    engine quirk, so we insert a throwaway character, wrap the bookmark
#    around it, then delete the character again - that leaves the bookmark
#    correctly collapsed at that position.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Create the Floor2.cs*") {
        $endPos = $p.Range.End - 1
        $tmp = $d.Range($endPos, $endPos)
        $tmp.InsertAfter("X")
        $tmpRange = $d.Range($endPos, $endPos + 1)
        $d.Bookmarks.Add("_GoBack", $tmpRange)
        $killRange = $d.Range($endPos, $endPos + 1)
        $killRange.Text = ""
    }
}

# ---------------------------------------------------------------------------
# 3) "Refactor the existing Unit code..." paragraph: the old bookmark split
#    the first sentence into two runs ("...You may " / "want to move...");
#    now that the bookmark has moved away, merge that text back into one
#    run and apply strike-through to the whole paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Refactor the existing Unit code*") {
        $startPos = $p.Range.Start
        $merged = "Refactor the existing Unit code. You may want to move some of the Unit functions up to the "
        $splitRange = $d.Range($startPos, $startPos + $merged.Length)
        $splitRange.Delete()
        $insertPoint = $d.Range($startPos, $startPos)
        $insertPoint.InsertBefore($merged)
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Refactor the existing Unit code*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 4-6) Mark the remaining finished to-do items with strike-through.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Perform tons of testing*" -or `
        $t -like "*Push changes*" -or `
        $t -like "*Delete the old walls, floors, units, and players*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 7) "Do a huge comment/class diagram sweep" is the new in-progress item -
#    give it the green highlight (moved from step 1's paragraph).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Do a huge comment*class diagram sweep*") {
        $p.Range.Font.HighlightColorIndex = 4   # wdBrightGreen -> w:val="green"
    }
}
